# Sync attendance_reports: reorder "Recorded By" (column G) values so that
# "System" is listed before the email/username entries instead of after.
#
# Rule (derived from the commit diff): for each comma-separated "Recorded By"
# cell, if it contains an exact-case "System" token, pull it out and place it
# first -- unless the first remaining token is the lowercase "system" token,
# in which case "System" is inserted right after it (so "system" stays the
# lead token). Cells without an exact "System" token are left untouched.

function Convert-RecordedBy($value) {
    if ($value -eq $null) {
        return $value
    }

    $parts = @($value -split ", ")

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $hasSystem = $true
        }
    }
    if (-not $hasSystem) {
        return $value
    }

    $rest = @()
    foreach ($p in $parts) {
        if (-not $p.Equals("System")) {
            $rest += $p
        }
    }

    if ($rest.Count -gt 0 -and $rest[0].Equals("system")) {
        if ($rest.Count -gt 1) {
            $newParts = @($rest[0], "System") + $rest[1..($rest.Count - 1)]
        } else {
            $newParts = @($rest[0], "System")
        }
    } else {
        $newParts = @("System") + $rest
    }

    return ($newParts -join ", ")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $current = $cell.Text
    if ([string]::IsNullOrEmpty($current)) {
        continue
    }
    $updated = Convert-RecordedBy $current
    if (-not $updated.Equals($current)) {
        $cell.Value = $updated
    }
}
